$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores plain-looking decimal numbers as literal text
# in the source data (e.g. "572.97", "0.120"). Excel auto-converts numeric-
# looking input to a true number (losing trailing zeros / introducing binary
# float noise), so force Text format on those cells before writing the value.
$textPriceCells = @("D5", "D6", "D7", "D12", "D14", "D16", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D35", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D48", "D51")
foreach ($cellRef in $textPriceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "68.821.02"
$ws.Range("E2").Value = "  -2.32%  "
$ws.Range("D3").Value = "3.440.16"
$ws.Range("E3").Value = "  -5.18%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "572.97"
$ws.Range("E5").Value = "  -4.92%  "
$ws.Range("D6").Value = "189.21"
$ws.Range("E6").Value = "  -3.95%  "
$ws.Range("D7").Value = "0.607"
$ws.Range("E7").Value = "  -2.99%  "
$ws.Range("D8").Value = "3.430.08"
$ws.Range("E8").Value = "  -5.17%  "
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("E10").Value = "  -7.25%  "
$ws.Range("E11").Value = "  -5.08%  "
$ws.Range("D12").Value = "50.80"
$ws.Range("E12").Value = "  -4.91%  "
$ws.Range("E13").Value = "  -7.82%  "
$ws.Range("D14").Value = "9.01"
$ws.Range("E14").Value = "  -5.47%  "
$ws.Range("D15").Value = "3.980.07"
$ws.Range("E15").Value = "  -5.18%  "
$ws.Range("D16").Value = "636.04"
$ws.Range("E16").Value = "  +4.87%  "
$ws.Range("D17").Value = "68.625.67"
$ws.Range("E17").Value = "  -2.71%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.433.73"
$ws.Range("E18").Value = "  -5.74%  "
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").Value = "0.120"
$ws.Range("E19").Value = "  -2.25%  "
$ws.Range("D20").Value = "12.14"
$ws.Range("E20").Value = "  -6.57%  "
$ws.Range("D21").Value = "17.95"
$ws.Range("E21").Value = "  -5.45%  "
$ws.Range("D22").Value = "0.934"
$ws.Range("E22").Value = "  -6.25%  "
$ws.Range("D23").Value = "17.62"
$ws.Range("E23").Value = "  -2.76%  "
$ws.Range("D24").Value = "5.28"
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("D25").Value = "97.86"
$ws.Range("E25").Value = "  -5.51%  "
$ws.Range("D26").Value = "4.23"
$ws.Range("E26").Value = "  -8.17%  "
$ws.Range("B27").Value = "ImmutableX"
$ws.Range("C27").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D27").Value = "2.81"
$ws.Range("E27").Value = "  -5.86%  "
$ws.Range("B28").Value = "LEO"
$ws.Range("C28").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D28").Value = "6.07"
$ws.Range("E28").Value = "  +1.88%  "
$ws.Range("D29").Value = "9.78"
$ws.Range("E29").Value = "  -7.77%  "
$ws.Range("D30").Value = "9.12"
$ws.Range("E30").Value = "  -5.97%  "
$ws.Range("D31").Value = "32.08"
$ws.Range("E31").Value = "  -5.01%  "
$ws.Range("D32").Value = "4.17"
$ws.Range("E32").Value = "  -11.01%  "
$ws.Range("D33").Value = "6.62"
$ws.Range("E33").Value = "  -8.04%  "
$ws.Range("E34").Value = "  -6.61%  "
$ws.Range("D35").Value = "60.73"
$ws.Range("E36").Value = "  -6.89%  "
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("D38").Value = "3.640.05"
$ws.Range("E38").Value = "  -8.00%  "
$ws.Range("D39").Value = "0.0₃0774"
$ws.Range("E39").Value = "  -12.21%  "
$ws.Range("D40").Value = "486.97"
$ws.Range("E40").Value = "  -5.77%  "
$ws.Range("D41").Value = "2.85"
$ws.Range("E41").Value = "  -7.52%  "
$ws.Range("D42").Value = "3.42"
$ws.Range("E42").Value = "  -3.20%  "
$ws.Range("B43").Value = "CoreDAO"
$ws.Range("C43").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D43").Value = "3.53"
$ws.Range("E43").Value = "  +69.42%  "
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").Value = "0.364"
$ws.Range("E44").Value = "  -6.48%  "
$ws.Range("B45").Value = "Kaspa"
$ws.Range("C45").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D45").Value = "0.132"
$ws.Range("E45").Value = "  -3.09%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Value = "33.89"
$ws.Range("E46").Value = "  -7.30%  "
$ws.Range("E47").Value = "  -5.84%  "
$ws.Range("D48").Value = "3.32"
$ws.Range("E48").Value = "  -5.58%  "
$ws.Range("E49").Value = "  -4.61%  "
$ws.Range("E50").Value = "  -5.05%  "
$ws.Range("D51").Value = "0.998"
$ws.Range("E51").Value = "  -0.48%  "
